# issue #5: stock data output to json file
#
# The "股票" (stock) sheet gains a new "property_category" column (constant
# value "stock" for every data row), inserted right after the existing
# "total" column and before "date" — pushing date / legislator_name /
# legislator_id one column to the right. Two company-name values also get
# a stray embedded space removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column H (shifting date/legislator_name/legislator_id -> I/J/K)
$ws.Columns.Item(8).Insert()

# New header + constant value for every data row (rows 2-6)
$ws.Range("H1").Value = "property_category"
$ws.Range("H2:H6").Value = "stock"

# Clean up two company names that had a stray embedded space
$ws.Range("B2").Value = "中美矽晶製品股份有限公司"
$ws.Range("B3").Value = "聯成化學科技股份有限公司"
